$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update row 4 (was the AF484509/Uganda/1998 subtype-A row) ---
# Now becomes a Brazilian subtype-C isolate. Accession written first.
$ws.Range("A4").Value = "U52953"

# --- Update row 3 (was the NC_001802/France/1983 duplicate row) ---
# Now becomes a Kenyan subtype-A isolate.
$ws.Range("A3").Value = "AF004885"
$ws.Range("E3").Value = "A"
$ws.Range("F3").Value = "Q23-CxC"
$ws.Range("G3").Value = 1994
$ws.Range("H3").Value = "Kenya"
$ws.Range("K3").Value = 12487816

# --- Finish row 4 ---
$ws.Range("E4").Value = "C"
$ws.Range("H4").Value = "Brazil"
$ws.Range("F4").Value = "92BR025"
$ws.Range("G4").Value = 1992
$ws.Range("K4").Value = 8891112

# --- Remove row 5 entirely (U46016 / Ethiopia / subtype C) ---
# Remaining rows below shift up by one.
$ws.Rows.Item(5).Delete()

# --- Narrow the scope: widen column A now that it holds fewer/longer entries ---
# (18.14 round-trips through Excel's pixel-quantised column width model to the
# stored width of 19 character-units.)
$ws.Columns.Item(1).ColumnWidth = 18.14

# --- Update selection to match the new, smaller used range ---
$ws.Range("A1:K9").Select()
